$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: add the "?" marker in column B (matches rows 3-5 pattern)
$ws.Range("B6").Value = "?"

# Row 7 (new user 5 results)
$ws.Range("B7").Value = "?"
$ws.Range("C7").Value = 39
$ws.Range("D7").Value = 318
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 195

# Reflect final active selection as in the authored workbook
$ws.Range("B8").Select()
